# Apply the updated cryptocurrency price/volume snapshot values.
# Source data is text (coin name / link / price / % volume columns),
# stored as inline strings in the sheet -- use .Value so Excel keeps
# them as text. For price cells that look like plain decimal numbers
# (e.g. "582.19", "0.0430"), force the cell to Text format first so
# Excel doesn't silently coerce the string into a Double (which would
# both lose trailing zeros and introduce floating point noise), then
# restore the cell's default style so no stray formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.440.76"
$ws.Range("E2").Value = "  +5.26%  "

# Row 3
$ws.Range("D3").Value = "3.479.51"
$ws.Range("E3").Value = "  +4.94%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.35%  "

# Row 7
$ws.Range("E7").Value = "  +2.88%  "

# Row 8
$ws.Range("D8").Value = "3.474.39"
$ws.Range("E8").Value = "  +5.01%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.16%  "

# Row 13
$ws.Range("E13").Value = "  -0.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.18%  "

# Row 15
$ws.Range("D15").Value = "4.014.78"
$ws.Range("E15").Value = "  +4.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.79%  "

# Row 17
$ws.Range("D17").Value = "3.468.60"
$ws.Range("E17").Value = "  +5.33%  "

# Row 18
$ws.Range("D18").Value = "67.263.29"
$ws.Range("E18").Value = "  +5.47%  "

# Row 19
$ws.Range("E19").Value = "  +0.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.61%  "

# Row 21
$ws.Range("E21").Value = "  +4.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "482.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.89%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.31%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +23.74%  "

# Row 25
$ws.Range("E25").Value = "  +9.67%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.00%  "

# Row 28
$ws.Range("E28").Value = "  +3.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.44%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.41%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.73%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "594.94"
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "  +5.51%  "

# Row 36
$ws.Range("E36").Value = "  +6.22%  "

# Row 37
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.55%  "

# Row 39
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.388"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.90%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.76%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0773"
$ws.Range("E41").Value = "  +6.98%  "

# Row 42
$ws.Range("D42").Value = "3.243.30"
$ws.Range("E42").Value = "  +6.90%  "

# Row 43
$ws.Range("E43").Value = "  +6.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0430"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.62%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.69%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +23.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "

# Row 48
$ws.Range("E48").Value = "  +2.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.18%  "

# Row 50
$ws.Range("E50").Value = "  +12.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "

